$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param(
        $Sheet,
        [string]$Addr,
        [string]$Text
    )
    $range = $Sheet.Range($Addr)
    $range.NumberFormat = "@"
    $range.Value = $Text
}

Set-CellText $ws "D2" "62.989.92"
Set-CellText $ws "E2" "  -0.07%  "
Set-CellText $ws "D3" "2.592.65"
Set-CellText $ws "E3" "  +1.67%  "
Set-CellText $ws "E4" "  -0.06%  "
Set-CellText $ws "D5" "585.11"
Set-CellText $ws "E5" "  -0.09%  "
Set-CellText $ws "D6" "147.19"
Set-CellText $ws "E6" "  +0.03%  "
Set-CellText $ws "D7" "0.999"
Set-CellText $ws "E7" "  -0.05%  "
Set-CellText $ws "E8" "  +2.30%  "
Set-CellText $ws "E9" "  +1.99%  "
Set-CellText $ws "E10" "  +2.05%  "
Set-CellText $ws "E11" "  -0.04%  "
Set-CellText $ws "E12" "  -0.07%  "
Set-CellText $ws "D13" "27.32"
Set-CellText $ws "E13" "  -0.26%  "
Set-CellText $ws "D14" "3.055.86"
Set-CellText $ws "E14" "  +1.66%  "
Set-CellText $ws "D15" "62.875.73"
Set-CellText $ws "E15" "  -0.09%  "
Set-CellText $ws "E16" "  +2.86%  "
Set-CellText $ws "D17" "2.592.99"
Set-CellText $ws "E17" "  +1.62%  "
Set-CellText $ws "D18" "11.32"
Set-CellText $ws "E18" "  -0.23%  "
Set-CellText $ws "D19" "342.36"
Set-CellText $ws "E19" "  +1.80%  "
Set-CellText $ws "D20" "4.40"
Set-CellText $ws "E20" "  +1.72%  "
Set-CellText $ws "D21" "6.70"
Set-CellText $ws "E21" "  -1.05%  "
Set-CellText $ws "E22" "  -0.02%  "
Set-CellText $ws "E23" "  +2.05%  "
Set-CellText $ws "D24" "2.714.15"
Set-CellText $ws "E24" "  +1.80%  "
Set-CellText $ws "E25" "  -1.50%  "
Set-CellText $ws "E26" "  -1.47%  "
Set-CellText $ws "E27" "  -0.03%  "
Set-CellText $ws "D28" "8.35"
Set-CellText $ws "E28" "  -0.23%  "
Set-CellText $ws "D29" "7.87"
Set-CellText $ws "E29" "  +5.14%  "
Set-CellText $ws "E30" "  -1.85%  "
Set-CellText $ws "E31" "  +0.46%  "
Set-CellText $ws "D32" "477.42"
Set-CellText $ws "E32" "  +15.13%  "
Set-CellText $ws "D33" "0.0₃0825"
Set-CellText $ws "E33" "  +1.51%  "
Set-CellText $ws "D34" "176.88"
Set-CellText $ws "E34" "  -0.56%  "
Set-CellText $ws "E35" "  +4.83%  "
Set-CellText $ws "E36" "  +0.07%  "
Set-CellText $ws "E37" "  +1.18%  "
Set-CellText $ws "D38" "19.05"
Set-CellText $ws "E38" "  -0.45%  "
Set-CellText $ws "D39" "4.53"
Set-CellText $ws "E39" "  +4.35%  "
Set-CellText $ws "D41" "1.71"
Set-CellText $ws "E41" "  -2.04%  "
Set-CellText $ws "D42" "158.72"
Set-CellText $ws "E42" "  +5.16%  "
Set-CellText $ws "E43" "  -0.01%  "
Set-CellText $ws "D44" "21.35"
Set-CellText $ws "E44" "  +2.40%  "
Set-CellText $ws "D45" "0.634"
Set-CellText $ws "E45" "  +5.45%  "
Set-CellText $ws "E46" "  +0.52%  "
Set-CellText $ws "E47" "  -0.07%  "
Set-CellText $ws "E48" "  -0.90%  "
Set-CellText $ws "D49" "18.39"
Set-CellText $ws "E49" "  +0.62%  "
Set-CellText $ws "E50" "  +1.19%  "
Set-CellText $ws "E51" "  +1.04%  "
